# This workbook's "Sold Qty" report lists the same SKU multiple times
# (different batches/rates). The edit re-orders which batch-row gets which
# Batch/Rate/Qty/Value figures by swapping the B:G (Batch No, MRP, Rate,
# Qty, Value) contents between specific row pairs - the A (Sr. No.) and
# H:M columns stay put since those describe the physical row position.
#
# For the 161/162/163 trio the rows rotate (161<-162<-163<-161), which we
# get by chaining two adjacent pairwise swaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($sheet, $row1, $row2) {
    $rng1 = $sheet.Range("B$row1`:G$row1")
    $rng2 = $sheet.Range("B$row2`:G$row2")
    $val1 = $rng1.Value2
    $val2 = $rng2.Value2
    $rng1.Value = $val2
    $rng2.Value = $val1
}

$pairs = @(
    @(149,150),
    @(279,280),
    @(346,347),
    @(350,352),
    @(355,356),
    @(375,376),
    @(379,380),
    @(389,390),
    @(400,401),
    @(419,420),
    @(431,432),
    @(457,458),
    @(536,537),
    @(590,591),
    @(599,600),
    @(601,602)
)

foreach ($pair in $pairs) {
    Swap-Rows $ws $pair[0] $pair[1]
}

# Rows 161-162-163 rotate: new161=old162, new162=old163, new163=old161.
# Two chained adjacent swaps give exactly that 3-cycle.
Swap-Rows $ws 161 162
Swap-Rows $ws 162 163
